$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at J:K (shifts existing PLZ..HASSH right to L..P)
$ws.Range("J1:K1").EntireColumn.Insert()

# Headers
$ws.Range("K1").Value = "HAUSNR"
$ws.Range("J1").Value = "STRASSE"

# Row 2 - Ackerstrasse 11
$ws.Range("J2").Value = "Ackerstrasse"
$ws.Range("K2").Value = 11

# Row 3 - Bertastrasse 22
$ws.Range("J3").Value = "Bertastrasse"
$ws.Range("K3").Value = 22

# Row 4 - Clausiensteig 3c (HAUSNR is text)
$ws.Range("J4").Value = "Clausiensteig"
$ws.Range("K4").Value = "3c"

# Row 6 - Dammweg 4
$ws.Range("J6").Value = "Dammweg"
$ws.Range("K6").Value = 4

# Row 7 - Erismannstrasse 505
$ws.Range("J7").Value = "Erismannstrasse"
$ws.Range("K7").Value = 505

# Row 8 - Floragasse 6 f (HAUSNR is text)
$ws.Range("J8").Value = "Floragasse"
$ws.Range("K8").Value = "6 f"

# Column widths for J:K (matches the width used for column I)
$ws.Range("J1:K1").EntireColumn.ColumnWidth = 22.33

# Selection
$ws.Range("K12").Select()
